$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 136.92308
$ws.Range("I6").Value = 136.92308
$ws.Range("K6").Value = 410.76924
$ws.Range("M6").Value = -298.76924
$ws.Range("H51").Value = 16998.75
$ws.Range("J51").Value = 9994.333000000001
$ws.Range("L51").Value = 9994.333000000001
$ws.Range("N51").Value = -10962.333
$ws.Range("H80").Value = 2784.4
$ws.Range("I80").Value = 336.42856
$ws.Range("K80").Value = 1009.28568
$ws.Range("M80").Value = -11.28567999999996
$ws.Range("H83").Value = 2784.4
$ws.Range("I83").Value = 336.42856
$ws.Range("K83").Value = 3027.85704
$ws.Range("M83").Value = 1964.14296
$ws.Range("H111").Value = 1404.7646
$ws.Range("I111").Value = 1106.75
$ws.Range("K111").Value = 3320.25
$ws.Range("M111").Value = -253.25
$ws.Range("H132").Value = 2881.7
$ws.Range("I132").Value = 2823.889
$ws.Range("K132").Value = 8471.667000000001
$ws.Range("M132").Value = -5941.667000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2288.182
$ws.Range("I32").Value = 2501.55
$ws.Range("J32").Value = 154.5
$ws.Range("K32").Value = 2501.55
$ws.Range("L32").Value = 154.5
$ws.Range("M32").Value = -2214.55
$ws.Range("N32").Value = -728.5
$ws.Range("H45").Value = 2180.6667
$ws.Range("J45").Value = 2180.2
$ws.Range("L45").Value = 2180.2
$ws.Range("N45").Value = -2934.2
$ws.Range("H74").Value = 5788462
$ws.Range("I74").Value = 3088036.8
$ws.Range("K74").Value = 3088036.8
$ws.Range("M74").Value = -3087162.8
$ws.Range("H77").Value = 5788462
$ws.Range("I77").Value = 3088036.8
$ws.Range("K77").Value = 15440184
$ws.Range("M77").Value = -15435816
$ws.Range("H97").Value = 542.1818
$ws.Range("J97").Value = 231
$ws.Range("L97").Value = 231
$ws.Range("N97").Value = -1223
$ws.Range("H122").Value = 2056.2856
$ws.Range("I122").Value = 2097.3333
$ws.Range("J122").Value = 1810
$ws.Range("K122").Value = 6291.999899999999
$ws.Range("L122").Value = 5430
$ws.Range("M122").Value = -3841.999899999999
$ws.Range("N122").Value = -10330
$ws.Range("H132").Value = 15155017
$ws.Range("I132").Value = 2920.9678
$ws.Range("K132").Value = 8762.903399999999
$ws.Range("M132").Value = -6232.903399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 481.4
$ws.Range("I22").Value = 326.66666
$ws.Range("K22").Value = 326.66666
$ws.Range("M22").Value = -153.66666
$ws.Range("H134").Value = 26522172
$ws.Range("I134").Value = 12507443
$ws.Range("K134").Value = 37522329
$ws.Range("M134").Value = -37519794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4145.2856
$ws.Range("I22").Value = 3122.5
$ws.Range("J22").Value = 4554.4
$ws.Range("K22").Value = 3122.5
$ws.Range("L22").Value = 4554.4
$ws.Range("M22").Value = -2772.5
$ws.Range("N22").Value = -5254.4
$ws.Range("H58").Value = 2431.4
$ws.Range("I58").Value = 2267.1333
$ws.Range("J58").Value = 2924.2
$ws.Range("K58").Value = 2267.1333
$ws.Range("L58").Value = 2924.2
$ws.Range("M58").Value = -2064.1333
$ws.Range("N58").Value = -3330.2
$ws.Range("H136").Value = 2431.4
$ws.Range("I136").Value = 2267.1333
$ws.Range("J136").Value = 2924.2
$ws.Range("K136").Value = 6801.3999
$ws.Range("L136").Value = 8772.599999999999
$ws.Range("M136").Value = -4251.3999
$ws.Range("N136").Value = -13872.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2102049
$ws.Range("I121").Value = 126311.375
$ws.Range("J121").Value = 10005000
$ws.Range("K121").Value = 378934.125
$ws.Range("L121").Value = 30015000
$ws.Range("M121").Value = -377624.125
$ws.Range("N121").Value = -30017620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2581.7273
$ws.Range("I80").Value = 2650.1667
$ws.Range("J80").Value = 2499.6
$ws.Range("K80").Value = 2650.1667
$ws.Range("L80").Value = 2499.6
$ws.Range("M80").Value = -1652.1667
$ws.Range("N80").Value = -4495.6
$ws.Range("H83").Value = 2581.7273
$ws.Range("I83").Value = 2650.1667
$ws.Range("J83").Value = 2499.6
$ws.Range("K83").Value = 13250.8335
$ws.Range("L83").Value = 12498
$ws.Range("M83").Value = -8258.833500000001
$ws.Range("N83").Value = -22482
$ws.Range("H97").Value = 603.8182
$ws.Range("I97").Value = 720.7143
$ws.Range("J97").Value = 399.25
$ws.Range("K97").Value = 720.7143
$ws.Range("L97").Value = 399.25
$ws.Range("M97").Value = -224.7143
$ws.Range("N97").Value = -1391.25
$ws.Range("H122").Value = 2175.6287
$ws.Range("I122").Value = 2373.2917
$ws.Range("J122").Value = 1744.3636
$ws.Range("K122").Value = 7119.875100000001
$ws.Range("L122").Value = 5233.0908
$ws.Range("M122").Value = -4669.875100000001
$ws.Range("N122").Value = -10133.0908
$ws.Range("H132").Value = 2252.6743
$ws.Range("I132").Value = 2215.8
$ws.Range("J132").Value = 2337.7693
$ws.Range("K132").Value = 6647.400000000001
$ws.Range("L132").Value = 7013.3079
$ws.Range("M132").Value = -4117.400000000001
$ws.Range("N132").Value = -12073.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 4143.9165
$ws.Range("J9").Value = 5247.5557
$ws.Range("L9").Value = 5247.5557
$ws.Range("N9").Value = -5695.5557
$ws.Range("H82").Value = 3836.2144
$ws.Range("I82").Value = 3240.4
$ws.Range("K82").Value = 3240.4
$ws.Range("M82").Value = -2879.4
$ws.Range("H85").Value = 3836.2144
$ws.Range("I85").Value = 3240.4
$ws.Range("K85").Value = 3240.4
$ws.Range("M85").Value = -1992.4
$ws.Range("H93").Value = 1596.56
$ws.Range("I93").Value = 1510
$ws.Range("J93").Value = 1645.25
$ws.Range("K93").Value = 1510
$ws.Range("L93").Value = 1645.25
$ws.Range("M93").Value = -262
$ws.Range("N93").Value = -4141.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7652.421
$ws.Range("I81").Value = 8850.571
$ws.Range("J81").Value = 4297.6
$ws.Range("K81").Value = 17701.142
$ws.Range("L81").Value = 8595.200000000001
$ws.Range("M81").Value = -16640.142
$ws.Range("N81").Value = -10717.2
$ws.Range("H84").Value = 7652.421
$ws.Range("I84").Value = 8850.571
$ws.Range("J84").Value = 4297.6
$ws.Range("K84").Value = 88505.70999999999
$ws.Range("L84").Value = 42976
$ws.Range("M84").Value = -83201.70999999999
$ws.Range("N84").Value = -53584
$ws.Range("H122").Value = 6965.6665
$ws.Range("I122").Value = 6965.6665
$ws.Range("K122").Value = 20896.9995
$ws.Range("M122").Value = -18446.9995
$ws.Range("H132").Value = 1377.5428
$ws.Range("I132").Value = 1103.4193
$ws.Range("K132").Value = 3310.2579
$ws.Range("M132").Value = -780.2579000000001
